{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph that currently reads the \"pledged to the campaign\"\n// sentence (it is split into multiple runs around the spell-checked word\n// \"company\") and the paragraph right after it that needs to be removed\n// entirely, plus one surplus blank paragraph that follows.\nlet targetIndex = -1;\nlet removeIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (targetIndex === -1 && t.indexOf(\"pledged to the campaign.\") !== -1) {\n    targetIndex = i;\n  } else if (\n    targetIndex !== -1 &&\n    removeIndex === -1 &&\n    t.indexOf(\"We can also assess if large fundings are successful or small funds.\") !== -1\n  ) {\n    removeIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1 || removeIndex === -1) {\n  throw new Error(\"Could not locate expected paragraphs.\");\n}\n\n// Re-insert the sentence as a single run (removes the spell-check run\n// splitting / proofErr markers around \"company\").\nconst targetParagraph = paras.items[targetIndex];\ntargetParagraph.insertText(\n  \"We don\\u2019t know how much money each individual or company pledged to the campaign.\",\n  Word.InsertLocation.replace\n);\n\n// Delete the whole \"We can also assess...\" paragraph.\nconst removeParagraph = paras.items[removeIndex];\nremoveParagraph.delete();\n\n// Delete one of the (now three) blank paragraphs that directly follow it,\n// bringing the run of blank paragraphs down to two.\nconst blankParagraph = paras.items[removeIndex + 1];\nblankParagraph.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that currently reads the \"pledged to the campaign\"\n# sentence (it is split into multiple runs around the spell-checked word\n# \"company\") and the paragraph right after it that needs to be removed\n# entirely, by scanning for their text (robust to any index drift).\n$targetIndex = -1\n$removeIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($targetIndex -eq -1 -and $t.Contains(\"pledged to the campaign.\")) {\n        $targetIndex = $i\n    } elseif ($targetIndex -ne -1 -and $removeIndex -eq -1 -and $t.Contains(\"We can also assess if large fundings are successful or small funds.\")) {\n        $removeIndex = $i\n        break\n    }\n}\n\n# Collapse the split runs of the target paragraph into a single run (removes\n# the spell-check run splitting / proofErr markers around \"company\"). The\n# visible text is already correct, so force a real mutation by writing a\n# placeholder first, then the desired final text.\n$p = $d.Paragraphs.Item($targetIndex)\n$r = $p.Range\n$full = $d.Range($r.Start, $r.End - 1)\n$placeholder = \"TEMP_PLACEHOLDER_TEXT\"\n$full.Text = $placeholder\n$finalText = \"We don\" + [char]8217 + \"t know how much money each individual or company pledged to the campaign.\"\n$finalRange = $d.Range($r.Start, $r.Start + $placeholder.Length)\n$finalRange.Text = $finalText\n\n# Delete the whole \"We can also assess...\" paragraph.\n$d.Paragraphs.Item($removeIndex).Range.Delete()\n\n# Delete one of the (now three) blank paragraphs that directly follow it,\n# bringing the run of blank paragraphs down to two.\n$d.Paragraphs.Item($removeIndex).Range.Delete()\n"}
